# Commit: "support pageable in each directive"
# - Sets A1 (the cell carrying the jx:area directive comment) to a test
#   marker value "ddddd".
# - Moves the active selection to B4:E4 (the merged HEADCOUNT value cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ddddd"
$ws.Range("B4:E4").Select()
